# Refresh the cryptocurrency price / volume(1h) snapshot table (automated data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.757.31'
$ws.Range('E2').Value = '  -4.13%  '
$ws.Range('D3').Value = '1.813.07'
$ws.Range('E3').Value = '  -3.20%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '276.39'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -8.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5048'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -5.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3505'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -6.57%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.35'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06671'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -6.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.94'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -7.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.8304'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -6.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07895'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.48%  '
$ws.Range('D14').Value = '1.821.88'
$ws.Range('E14').Value = '  -2.70%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.070'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.87%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '87.46'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -6.06%  '
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.01'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -4.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008040'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -5.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').Value = '25.802.97'
$ws.Range('E21').Value = '  -4.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.717'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -5.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.983'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -6.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.056'
$ws.Range('D24').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.85'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.168'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -4.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.667'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.49%  '
$ws.Range('E28').Value = '  -5.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '109.58'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -4.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.325'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -8.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.224'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -7.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08808'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04855'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7258'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -9.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.131'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.64%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.874'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -3.90%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.150'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.71%  '
$ws.Range('B38').Value = 'Frax'
$ws.Range('C38').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.0000'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5198'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -11.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01842'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -5.61%  '
$ws.Range('E41').Value = '  -13.62%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9510'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -10.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '113.27'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.157'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -7.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.067'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -8.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9999'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4549'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -9.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1361'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -8.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.309'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -6.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.35'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.501'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -6.94%  '
